$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to remain text (matches original inlineStr / string-typed cells)
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
}

Set-TextValue $ws.Range("D2") "66.624.01"
$ws.Range("E2").Value = "  +0.58%  "
Set-TextValue $ws.Range("D3") "3.587.67"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue $ws.Range("D5") "609.48"
$ws.Range("E5").Value = "  +0.54%  "
Set-TextValue $ws.Range("D6") "146.49"
$ws.Range("E6").Value = "  +1.17%  "
Set-TextValue $ws.Range("D7") "3.587.11"
$ws.Range("E7").Value = "  +0.77%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("E10").Value = "  -0.02%  "
Set-TextValue $ws.Range("D11") "7.97"
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("E12").Value = "  +1.11%  "
Set-TextValue $ws.Range("D13") "4.194.51"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("E14").Value = "  +0.51%  "
Set-TextValue $ws.Range("D15") "30.07"
$ws.Range("E15").Value = "  -0.59%  "
Set-TextValue $ws.Range("D16") "3.587.90"
$ws.Range("E16").Value = "  +0.85%  "
Set-TextValue $ws.Range("D17") "66.688.51"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("E18").Value = "  +0.38%  "
Set-TextValue $ws.Range("D19") "11.45"
$ws.Range("E19").Value = "  -2.12%  "
Set-TextValue $ws.Range("D20") "6.26"
$ws.Range("E20").Value = "  +0.58%  "
Set-TextValue $ws.Range("D21") "14.98"
$ws.Range("E21").Value = "  +0.21%  "
Set-TextValue $ws.Range("D22") "433.05"
$ws.Range("E23").Value = "  +2.20%  "
Set-TextValue $ws.Range("D24") "79.11"
Set-TextValue $ws.Range("D25") "3.732.42"
$ws.Range("E25").Value = "  +0.91%  "
Set-TextValue $ws.Range("D26") "1.00"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  +0.34%  "
Set-TextValue $ws.Range("D28") "9.32"
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("E31").Value = "  +0.10%  "
Set-TextValue $ws.Range("D32") "3.583.32"
$ws.Range("E32").Value = "  +0.85%  "
Set-TextValue $ws.Range("D33") "25.51"
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("E34").Value = "  -3.34%  "
$ws.Range("E35").Value = "  -1.73%  "
Set-TextValue $ws.Range("D36") "7.85"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("E37").Value = "  +0.03%  "
Set-TextValue $ws.Range("D38") "1.72"
$ws.Range("E38").Value = "  -2.14%  "
Set-TextValue $ws.Range("D39") "5.64"
$ws.Range("E39").Value = "  -0.26%  "
Set-TextValue $ws.Range("D40") "174.13"
$ws.Range("E40").Value = "  +1.19%  "
Set-TextValue $ws.Range("D41") "0.0854"
$ws.Range("E41").Value = "  -0.35%  "
Set-TextValue $ws.Range("D42") "5.23"
$ws.Range("E42").Value = "  -1.39%  "
Set-TextValue $ws.Range("D43") "0.895"
$ws.Range("E43").Value = "  +0.00%  "
Set-TextValue $ws.Range("D44") "1.92"
$ws.Range("E44").Value = "  +0.77%  "
Set-TextValue $ws.Range("D45") "45.81"
$ws.Range("E45").Value = "  -0.11%  "
Set-TextValue $ws.Range("D46") "0.999"
$ws.Range("E46").Value = "  -0.02%  "
Set-TextValue $ws.Range("D47") "2.53"
$ws.Range("E47").Value = "  +5.52%  "
$ws.Range("E48").Value = "  -1.64%  "
Set-TextValue $ws.Range("D49") "24.98"
$ws.Range("E49").Value = "  -4.25%  "

# Rows 50 & 51: coin order swapped (EnergySwap now ranks above Cosmos), with updated values
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D50") "23.83"
$ws.Range("E50").Value = "  +4.13%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D51") "7.20"
$ws.Range("E51").Value = "  +0.81%  "
